$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.8733400000000001
$ws.Range("H2").Value = 2.62002
$ws.Range("I2").Value = 0.2319025556440181
$ws.Range("J2").Value = 0.2319025556440181
$ws.Range("M2").Value = 12.87437866666667
$ws.Range("N2").Value = 38.623136
$ws.Range("O2").Value = 0.272778495601419
$ws.Range("P2").Value = 0.272778495601419
$ws.Range("Q2").Value = 11.24370986474667
$ws.Range("R2").Value = 101.19338878272
$ws.Range("S2").Value = 0.06325803025469962
$ws.Range("T2").Value = 0.06325803025469962

$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.8733400000000001
$ws.Range("H3").Value = 2.62002
$ws.Range("I3").Value = 0.2319025556440181
$ws.Range("J3").Value = 0.2319025556440181
$ws.Range("O3").Value = 0.3751152716914535
$ws.Range("P3").Value = 0.3751152716914535
$ws.Range("Q3").Value = 15.46194934258
$ws.Range("R3").Value = 139.15754408322
$ws.Range("S3").Value = 0.08699019016634828
$ws.Range("T3").Value = 0.08699019016634826

$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.8733400000000001
$ws.Range("H4").Value = 2.62002
$ws.Range("I4").Value = 0.2319025556440181
$ws.Range("J4").Value = 0.2319025556440181
$ws.Range("M4").Value = 8.077278666666666
$ws.Range("N4").Value = 24.231836
$ws.Range("O4").Value = 0.1711389714636405
$ws.Range("P4").Value = 0.1711389714636405
$ws.Range("Q4").Value = 7.054210550746667
$ws.Range("R4").Value = 63.48789495672001
$ws.Range("S4").Value = 0.03968756485270693
$ws.Range("T4").Value = 0.03968756485270692

$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.8733400000000001
$ws.Range("H5").Value = 2.62002
$ws.Range("I5").Value = 0.2319025556440181
$ws.Range("J5").Value = 0.2319025556440181
$ws.Range("M5").Value = 8.541146333333332
$ws.Range("N5").Value = 25.623439
$ws.Range("O5").Value = 0.180967261243487
$ws.Range("P5").Value = 0.180967261243487
$ws.Range("Q5").Value = 7.459324738753333
$ws.Range("R5").Value = 67.13392264878
$ws.Range("S5").Value = 0.04196677037026331
$ws.Range("T5").Value = 0.0419667703702633

$ws.Range("I6").Value = 0.07113291057171066
$ws.Range("J6").Value = 0.07113291057171067
$ws.Range("M6").Value = 12.87437866666667
$ws.Range("N6").Value = 38.623136
$ws.Range("O6").Value = 0.272778495601419
$ws.Range("P6").Value = 0.272778495601419
$ws.Range("Q6").Value = 3.44885292912
$ws.Range("R6").Value = 31.03967636208
$ws.Range("S6").Value = 0.0194035283335015
$ws.Range("T6").Value = 0.01940352833350151

$ws.Range("I7").Value = 0.07113291057171066
$ws.Range("J7").Value = 0.07113291057171067
$ws.Range("O7").Value = 0.3751152716914535
$ws.Range("P7").Value = 0.3751152716914535
$ws.Range("S7").Value = 0.02668304107531111
$ws.Range("T7").Value = 0.02668304107531111

$ws.Range("I8").Value = 0.07113291057171066
$ws.Range("J8").Value = 0.07113291057171067
$ws.Range("M8").Value = 8.077278666666666
$ws.Range("N8").Value = 24.231836
$ws.Range("O8").Value = 0.1711389714636405
$ws.Range("P8").Value = 0.1711389714636405
$ws.Range("Q8").Value = 2.16378179562
$ws.Range("R8").Value = 19.47403616058
$ws.Range("S8").Value = 0.01217361315245768
$ws.Range("T8").Value = 0.01217361315245768

$ws.Range("I9").Value = 0.07113291057171066
$ws.Range("J9").Value = 0.07113291057171067
$ws.Range("M9").Value = 8.541146333333332
$ws.Range("N9").Value = 25.623439
$ws.Range("O9").Value = 0.180967261243487
$ws.Range("P9").Value = 0.180967261243487
$ws.Range("Q9").Value = 2.288044985505
$ws.Range("R9").Value = 20.592404869545
$ws.Range("S9").Value = 0.01287272801044036
$ws.Range("T9").Value = 0.01287272801044036

$ws.Range("G10").Value = 0.06721833333333334
$ws.Range("H10").Value = 0.201655
$ws.Range("I10").Value = 0.01784883697773089
$ws.Range("J10").Value = 0.01784883697773089
$ws.Range("M10").Value = 12.87437866666667
$ws.Range("N10").Value = 38.623136
$ws.Range("O10").Value = 0.272778495601419
$ws.Range("P10").Value = 0.272778495601419
$ws.Range("Q10").Value = 0.8653942766755557
$ws.Range("R10").Value = 7.78854849008
$ws.Range("S10").Value = 0.004868778899020409
$ws.Range("T10").Value = 0.004868778899020409

$ws.Range("G11").Value = 0.06721833333333334
$ws.Range("H11").Value = 0.201655
$ws.Range("I11").Value = 0.01784883697773089
$ws.Range("J11").Value = 0.01784883697773089
$ws.Range("O11").Value = 0.3751152716914535
$ws.Range("P11").Value = 0.3751152716914535
$ws.Range("Q11").Value = 1.190059386828334
$ws.Range("R11").Value = 10.710534481455
$ws.Range("S11").Value = 0.006695371332277983
$ws.Range("T11").Value = 0.006695371332277982

$ws.Range("G12").Value = 0.06721833333333334
$ws.Range("H12").Value = 0.201655
$ws.Range("I12").Value = 0.01784883697773089
$ws.Range("J12").Value = 0.01784883697773089
$ws.Range("M12").Value = 8.077278666666666
$ws.Range("N12").Value = 24.231836
$ws.Range("O12").Value = 0.1711389714636405
$ws.Range("P12").Value = 0.1711389714636405
$ws.Range("Q12").Value = 0.5429412098422223
$ws.Range("R12").Value = 4.88647088858
$ws.Range("S12").Value = 0.003054631602191058
$ws.Range("T12").Value = 0.003054631602191057

$ws.Range("G13").Value = 0.06721833333333334
$ws.Range("H13").Value = 0.201655
$ws.Range("I13").Value = 0.01784883697773089
$ws.Range("J13").Value = 0.01784883697773089
$ws.Range("M13").Value = 8.541146333333332
$ws.Range("N13").Value = 25.623439
$ws.Range("O13").Value = 0.180967261243487
$ws.Range("P13").Value = 0.180967261243487
$ws.Range("Q13").Value = 0.5741216212827778
$ws.Range("R13").Value = 5.167094591544999
$ws.Range("S13").Value = 0.003230055144241436
$ws.Range("T13").Value = 0.003230055144241436

$ws.Range("G14").Value = 2.557535
$ws.Range("H14").Value = 7.672605
$ws.Range("I14").Value = 0.6791156968065403
$ws.Range("J14").Value = 0.6791156968065403
$ws.Range("M14").Value = 12.87437866666667
$ws.Range("N14").Value = 38.623136
$ws.Range("O14").Value = 0.272778495601419
$ws.Range("P14").Value = 0.272778495601419
$ws.Range("Q14").Value = 32.92667404325334
$ws.Range("R14").Value = 296.34006638928
$ws.Range("S14").Value = 0.1852481581141974
$ws.Range("T14").Value = 0.1852481581141974

$ws.Range("G15").Value = 2.557535
$ws.Range("H15").Value = 7.672605
$ws.Range("I15").Value = 0.6791156968065403
$ws.Range("J15").Value = 0.6791156968065403
$ws.Range("O15").Value = 0.3751152716914535
$ws.Range("P15").Value = 0.3751152716914535
$ws.Range("Q15").Value = 45.27958940604501
$ws.Range("R15").Value = 407.516304654405
$ws.Range("S15").Value = 0.2547466691175161
$ws.Range("T15").Value = 0.2547466691175161

$ws.Range("G16").Value = 2.557535
$ws.Range("H16").Value = 7.672605
$ws.Range("I16").Value = 0.6791156968065403
$ws.Range("J16").Value = 0.6791156968065403
$ws.Range("M16").Value = 8.077278666666666
$ws.Range("N16").Value = 24.231836
$ws.Range("O16").Value = 0.1711389714636405
$ws.Range("P16").Value = 0.1711389714636405
$ws.Range("Q16").Value = 20.65792289475333
$ws.Range("R16").Value = 185.92130605278
$ws.Range("S16").Value = 0.1162231618562848
$ws.Range("T16").Value = 0.1162231618562848

$ws.Range("G17").Value = 2.557535
$ws.Range("H17").Value = 7.672605
$ws.Range("I17").Value = 0.6791156968065403
$ws.Range("J17").Value = 0.6791156968065403
$ws.Range("M17").Value = 8.541146333333332
$ws.Range("N17").Value = 25.623439
$ws.Range("O17").Value = 0.180967261243487
$ws.Range("P17").Value = 0.180967261243487
$ws.Range("Q17").Value = 21.84428068762167
$ws.Range("R17").Value = 196.598526188595
$ws.Range("S17").Value = 0.1228977077185419
$ws.Range("T17").Value = 0.1228977077185419

